# Fix Training Data Issue (#48)
#
# The "BF" column stores a per-row date label. Every data row (2-31) had
# been stamped with the source file's own name ("6-21-2007-08") instead of
# the actual game date. Because of how the NBA stats were originally
# scraped, the date was off by one day relative to the ISO date that
# should have been recorded. Correct every BF2:BF31 cell to the proper
# ISO-8601 date string "2008-06-21".
#
# The literal text "2008-06-21" looks like a date, so a plain
# Range.Value assignment would be auto-converted into a date serial
# number by Excel's normal "smart" cell-entry parsing. Priming the cell
# with a Text number format before the write keeps it a literal string
# (matching the original inline-string "Date" column), and resetting the
# style back to Normal afterwards keeps the cell's formatting identical
# to how it started (no format was applied to this column before the
# edit, and none should be applied after it either).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "6-21-2007-08"
$newValue = "2008-06-21"

# Data rows are 2-31 (row 1 is the "Date" header); column BF holds the
# mis-stamped date label for every row in this sheet.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF" + $row)
    if ($cell.Text -eq $oldValue) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}
